# Add a new "community IPT" intervention row to the constants sheet.
# This adds program_timeperiod_community_ipt_round, a time-period parameter
# for the community IPT intervention (not dependent on patients started on
# treatment), inserted just above the existing program_rate_start_treatment
# row (row 38).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 38 - everything from row 38 down shifts to
# row 39+, and the new row inherits formatting from the row above (row 37).
$ws.Rows("38:38").Insert() | Out-Null

# Populate the new row with the new parameter name/value.
$ws.Range("A38").Value = "program_timeperiod_community_ipt_round"
$ws.Range("B38").Value = 1

# Reflect the new cursor position / selection left behind by the edit.
$ws.Range("A38").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
